# Activity_Productivity.xlsx - "more test cases for activity api"
#
# This script reproduces the authored changes:
#  - two new schema-filename cells (and their new shared strings) on the
#    "GET_user_sessions" and "GET_user_simultaneous" sheets
#  - removal of the now-unused "ExpectedResponseBody" (column I) on the
#    "GET_user_simultaneous" sheet
#  - updated sheet view selections/scroll positions on several sheets
#  - "GET_user_simultaneous" becomes the active tab/sheet

$wb = $excel.ActiveWorkbook

$wsActivity     = $wb.Worksheets.Item("GET_activity_v1_users")
$wsLastLogin    = $wb.Worksheets.Item("GET_last_login")
$wsSimultaneous = $wb.Worksheets.Item("GET_user_simultaneous")
$wsSessions     = $wb.Worksheets.Item("GET_user_sessions")

# ---------------------------------------------------------------------
# GET_user_sessions: add the schema validation file name for row 2
# (new shared string -> "getUserSessionsSchema.json")
# ---------------------------------------------------------------------
$wsSessions.Range("H2").Value = "getUserSessionsSchema.json"

# ---------------------------------------------------------------------
# GET_user_simultaneous: add schema validation file names for rows 2/3
# and drop the now unused "ExpectedResponseBody" column (I)
# ---------------------------------------------------------------------
$wsSimultaneous.Range("H2").Value = "getSimultaniousActivitySchema.json"
$wsSimultaneous.Range("H3").Value = "400error.json"
$wsSimultaneous.Columns("I:I").Delete()

# ---------------------------------------------------------------------
# Sheet view / selection updates
# ---------------------------------------------------------------------

# GET_activity_v1_users: scroll to E1, select H4 (no longer the active tab)
$wsActivity.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$wsActivity.Range("H4").Select()

# GET_last_login: scroll to D1, select F15
$wsLastLogin.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$wsLastLogin.Range("F15").Select()

# GET_user_sessions: scroll to E1, select H12
$wsSessions.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$wsSessions.Range("H12").Select()

# GET_user_simultaneous: scroll to E1, select I1:I1048576, and leave this
# sheet active/selected (matches activeTab=2 on the workbook)
$wsSimultaneous.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$wsSimultaneous.Range("I1:I1048576").Select()

Write-Output "edit complete"
